# Products workbook update:
#  - Add 4 new product rows (rows 3-6) to Sheet1 / Table2
#  - Wrap-text the Description column (all rows) and the Name column for the
#    "Petrol Washer" rows (2-4)
#  - Reformat the Price column to show 2 decimal places
#  - Resize columns C/D/F to fit the new content
#  - Leave the selection positioned under the new data (A7), matching the
#    state the workbook was saved in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3 - Honda GX160 acid-pump pressure washer
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Petrol Washer"
$ws.Range("C3").Value = "Honda GX160 + MC18 Domestic Acid Pump Viton Seals Carry Frame + 30m Reel"
$ws.Range("D3").Value = 1045
$ws.Range("E3").Value = "This versatile machine is designed to tackle any cleaning task with a selection of nozzles for different spray patterns and an adjustable pressure regulator, to help you clean those hard to reach areas!"
$ws.Range("F3").Value = "/static/PetrolPowered/2.jpg"
$ws.Range("G3").Value = 10

# ---------------------------------------------------------------------------
# Row 4 - Honda GX200 semi-industrial pressure washer
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Petrol Washer"
$ws.Range("C4").Value = "Maxflow Semi-Industrial Pressure Washer - Honda GX200 14 LPM Low Profile Frame"
$ws.Range("D4").Value = 915
$ws.Range("E4").Value = "Maxflow cold water pressure washer powered by a Honda GX200 petrol engine. Featuring a premium series Comet pump, this machine produces 14 litres per minute water flow at a pressure of 150 bar. Includes chemical pick up hose complete with filter, and adjustable pressure nozzle for chemical application"
$ws.Range("F4").Value = "/static/PetrolPowered/3.jpg"
$ws.Range("G4").Value = 10

# ---------------------------------------------------------------------------
# Row 5 - Yanmar diesel industrial pressure washer
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Diesel Washer"
$ws.Range("C5").Value = "Maxflow Industrial Pressure Washer - Yanmar L100-V Comet Pump 18 LPM Trolley Frame + Reel"
$ws.Range("D5").Value = 2625
$ws.Range("E5").Value = "Maxflow cold water pressure washer powered by a Yanmar L100 V spec diesel engine. Featuring a HTD type belt driven premium series Comet pump, this machine produces 18 litres per minute water flow at a pressure of 200 bar. Includes reel, chemical pick up hose complete with filter, and adjustable pressure nozzle for chemical application"
$ws.Range("F5").Value = "/static/DieselWasher/4.jpg"
$ws.Range("G5").Value = 10

# ---------------------------------------------------------------------------
# Row 6 - Yanmar diesel hot-wash industrial pressure washer
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Diesel Washer"
$ws.Range("C6").Value = "Maxflow Industrial Diesel Hot Pressure Washer - Yanmar L100-V 18 LPM Trolley Frame + Reel"
$ws.Range("D6").Value = 7795
$ws.Range("E6").Value = "Maxflow diesel engine-driven hot wash machine for industrial applications. Produces up to 18 litres per minute at a maximum 200 bar pressure."
$ws.Range("F6").Value = "/static/DieselWasher/5.jpg"
$ws.Range("G6").Value = 10

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Price column now shows pence (2 decimal places)
$priceFormat = """£""#,##0.00;[Red]\-""£""#,##0.00"
$ws.Range("D2:D6").NumberFormat = $priceFormat

# Wrap the Description column (header + every data row)
$ws.Range("E1:E6").WrapText = $true

# Wrap the Name column only for the Petrol Washer rows (2-4)
$ws.Range("C2:C4").WrapText = $true

# Row heights for the new rows (row 6 keeps the sheet default height)
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45

# Column widths to fit the new, longer content
$ws.Columns.Item(3).ColumnWidth = 112.5
$ws.Columns.Item(4).ColumnWidth = 8.333333333333334
$ws.Columns.Item(6).ColumnWidth = 24.666666666666668

# Leave the selection where it ended up after entering the new rows
$ws.Range("A7").Select()

Write-Output "Added 4 rows to Table2 and refreshed formatting"
